$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 82, pushing all
# subsequent rows (old 82..174) down by one (new 83..175).
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A82").Value = 7
$ws.Range("B82").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C82").Value = "Ñuble"
$ws.Range("D82").Value2 = 44494
$ws.Range("E82").Value = 16
$ws.Range("F82").Value = 100112008
$ws.Range("G82").Value = "Coliflor"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 200
$ws.Range("K82").Value = 700
$ws.Range("L82").Value = 800
$ws.Range("M82").Value = 750
$ws.Range("N82").Value = "$/unidad"
$ws.Range("O82").Value = "Región Metropolitana"
$ws.Range("P82").Value = 750
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"
